$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 28
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025-09-19"
$ws.Cells.Item($row, 2).Value = "21:21:24"
$ws.Cells.Item($row, 3).Value = "1.00 EUR = 1,780.3126"
